# Applies the RHPF -> renamed / content update described in the commit:
#   "Fixes hydrogen/CCS sector bugs including those in (#347); renames RHPF
#    as BRHPF and NBRHPF as RHPF to match EPS naming convention"
#
# Content changes:
#  - About sheet: A12 note text changes to reflect repeal-of-IRA-tax-credits wording
#  - RHPF sheet: "hydrocarbon partial oxidation" pathway renamed to
#    "thermochemical water splitting" (header F1 and row label A6)
#  - RHPF sheet data: electrolysis row (row2) and natural gas reforming row (row3)
#    populated with 0.05 / 0.95 fractions across all columns (B:H), and the
#    "electrolysis with guaranteed clean electricity" row (row7) zeroed out,
#    with G/H columns picking up the right-aligned numeric style used by B:F.
#  - Selections on both sheets updated to match the saved state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("RHPF")

# --- About sheet: update the note that used to describe "guaranteed clean
# electricity" electrolysis, now describing the repeal-of-IRA-tax-credits note.
$ws1.Range("A12").Value = "the BAU production shares, representing a repeal of IRA tax credits"

# --- RHPF sheet: rename the "hydrocarbon partial oxidation" pathway to
# "thermochemical water splitting" (used in the F1 header and the A6 row label).
$ws2.Range("F1").Value = "thermochemical water splitting"
$ws2.Range("A6").Value = "thermochemical water splitting"

# --- RHPF sheet: update the fraction grid.
# Row 2 = electrolysis -> 0.05 for every column B..H
$ws2.Range("B2:H2").Value = 0.05

# Row 3 = natural gas reforming -> 0.95 for every column B..H
$ws2.Range("B3:H3").Value = 0.95

# Row 7 = electrolysis with guaranteed clean electricity -> 0 for every column B..H
$ws2.Range("B7:H7").Value = 0

# Match the style used by columns B:F (right-aligned) for G2:H2, G3:H3, G7:H7,
# since those cells previously had the default (unstyled) numeric format.
$ws2.Range("G2:H2").HorizontalAlignment = -4152
$ws2.Range("G3:H3").HorizontalAlignment = -4152
$ws2.Range("G7:H7").HorizontalAlignment = -4152

# --- Restore the selections recorded in the saved workbook.
# (About stays the active/tabSelected sheet, so select it last.)
$ws2.Range("B2:H3").Select()
$ws1.Range("B13").Select()
